$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8; this shifts the existing rows 8-64 down to 9-65,
# preserving all of their data and formatting (including the date style on column D).
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new record.
$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "Vega Modelo de Temuco"
$ws.Range("C8").Value = "La Araucanía"
$ws.Range("D8").Value = 45061
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100107
$ws.Range("H8").Value = "Otros"
$ws.Range("I8").Value = 100107001
$ws.Range("J8").Value = "Caqui"
$ws.Range("K8").Value = "Fuyu"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 140
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("Q8").Value = "$/bandeja 15 kilos granel"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 1200
$ws.Range("T8").Value = 15

# Ensure the date style used elsewhere in column D is applied to the new cell.
$ws.Range("D8").NumberFormat = $ws.Range("D9").NumberFormat
